$wb = $excel.ActiveWorkbook

# "Forecast Comparison" sheet - numeric MyForecast values
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D2").Value = 27
$wsForecast.Range("D11").Value = 41

# "Summary" sheet - these look like numbers/dates but are stored as plain
# text strings in the workbook. Force Text number format before assigning
# so Excel doesn't auto-convert them into real numbers/dates, then clear
# the format change afterwards so no stray style is left on the cells.
$wsSummary = $wb.Worksheets.Item("Summary")
$summaryRange = $wsSummary.Range("B9:B15")
$summaryRange.NumberFormat = "@"

$wsSummary.Range("B9").Value = "709"
$wsSummary.Range("B10").Value = "360"
$wsSummary.Range("B11").Value = "190"
$wsSummary.Range("B12").Value = "56"
$wsSummary.Range("B13").Value = "2025-03-09"
$wsSummary.Range("B14").Value = "27"
$wsSummary.Range("B15").Value = "2025-01-26"

$summaryRange.ClearFormats()
